$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "incident" table (rows 46-56) is getting a new field "incident_severity"
# inserted right after "incident_status" (row 52). Insert a new row at 53,
# pushing everything below it (incident_status's former neighbours, and the
# whole incident_comment table) down by one.
$ws.Rows("53:53").Insert()

# Copy the formatting (styles/borders/fill) from the row above (52) onto the
# freshly inserted row so it keeps the same "inner row of a merged block"
# look instead of Excel's bare insert-row default style.
$ws.Range("B52:F52").Copy() | Out-Null
$ws.Range("B53:F53").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# incident_status's description gains the word "진행" (진행 상태 instead of 상태)
$ws.Range("F52").Value = "인시던트 진행 상태 (1:계획되지 않음 2:해결중 3:완료 4:비활성화)"

# Fill in the new incident_severity row
$ws.Range("C53").Value = "incident_severity"
$ws.Range("D53").Value = ""
$ws.Range("E53").Value = "v"
$ws.Range("F53").Value = "인시던트 심각도 ( 1:Critical 2:Major 3:Minor ) "

# Column F needs to widen to fit the longer descriptions now in the sheet.
$ws.Columns("F").ColumnWidth = 60

# Restore the selection state to where the author ended up after the edit.
$ws.Range("M49").Select()
